$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New shared string "exp fit" plus the new data block (columns T:Z, rows 9-27)
# mirroring the existing E[Vcm-1]/alpha[cm-1] exponential-fit block in L:R,
# with a second "E[Vm-1]/alpha[m-1]" block in T:X rows 22-27.

$ws.Range("T9").Value = 'exp fit'
$ws.Range("T9").ClearFormats()

$ws.Range("T10").Value = 'E [Vcm-1]'
$ws.Range("T10").Font.Bold = $true

$ws.Range("U10").Value = 'alpha [cm-1]'
$ws.Range("U10").NumberFormat = '0.000'

$ws.Range("W10").Value = 'A'
$ws.Range("W10").ClearFormats()

$ws.Range("X10").Value = 'B'
$ws.Range("X10").NumberFormat = '0.000'

$ws.Range("Y10").Value = 'i0'
$ws.Range("Y10").NumberFormat = '0.000'

$ws.Range("Z10").Value = 'U0'
$ws.Range("Z10").NumberFormat = '0.000'

$ws.Range("T11").Value = 80
$ws.Range("T11").NumberFormat = '0.000'

$ws.Range("U11").Value = 1.60208
$ws.Range("U11").ClearFormats()

$ws.Range("Y11").Value = 63.65086
$ws.Range("Y11").ClearFormats()

$ws.Range("T12").Value = 'errors'
$ws.Range("T12").NumberFormat = '0.000'

$ws.Range("U12").Value = 0.05722
$ws.Range("U12").ClearFormats()

$ws.Range("Y12").Value = 6.62342
$ws.Range("Y12").ClearFormats()

$ws.Range("T13").Value = 90
$ws.Range("T13").NumberFormat = '0.000'

$ws.Range("U13").Value = 1.85248
$ws.Range("U13").ClearFormats()

$ws.Range("Y13").Value = 45.29801
$ws.Range("Y13").ClearFormats()

$ws.Range("T14").Value = 'errors'
$ws.Range("T14").NumberFormat = '0.000'

$ws.Range("U14").Value = 0.03326
$ws.Range("U14").ClearFormats()

$ws.Range("Y14").Value = 2.79297
$ws.Range("Y14").ClearFormats()

$ws.Range("T15").Value = 100
$ws.Range("T15").NumberFormat = '0.000'

$ws.Range("U15").Value = 2.06715
$ws.Range("U15").ClearFormats()

$ws.Range("Y15").Value = 26.65003
$ws.Range("Y15").ClearFormats()

$ws.Range("T16").Value = 'errors'
$ws.Range("T16").NumberFormat = '0.000'

$ws.Range("U16").Value = 0.03395
$ws.Range("U16").ClearFormats()

$ws.Range("Y16").Value = 1.69932
$ws.Range("Y16").ClearFormats()

$ws.Range("T17").Value = 110
$ws.Range("T17").NumberFormat = '0.000'

$ws.Range("U17").Value = 2.28404
$ws.Range("U17").ClearFormats()

$ws.Range("Y17").Value = 19.91821
$ws.Range("Y17").ClearFormats()

$ws.Range("T18").Value = 'errors'
$ws.Range("T18").NumberFormat = '0.000'

$ws.Range("U18").Value = 0.03244
$ws.Range("U18").ClearFormats()

$ws.Range("Y18").Value = 1.22679
$ws.Range("Y18").ClearFormats()

$ws.Range("T19").Value = 120
$ws.Range("T19").NumberFormat = '0.000'

$ws.Range("U19").Value = 2.78306
$ws.Range("U19").ClearFormats()

$ws.Range("Y19").Value = 7.088
$ws.Range("Y19").ClearFormats()

$ws.Range("T20").Value = 'errors'
$ws.Range("T20").NumberFormat = '0.000'

$ws.Range("U20").Value = 0.03503
$ws.Range("U20").ClearFormats()

$ws.Range("Y20").Value = 0.48057
$ws.Range("Y20").ClearFormats()

$ws.Range("T21").NumberFormat = '0.000'

$ws.Range("T22").Value = 'E [Vm-1]'
$ws.Range("T22").Font.Bold = $true

$ws.Range("U22").Value = 'alpha [m-1]'
$ws.Range("U22").NumberFormat = '0.000'

$ws.Range("V22").Value = 'alpha/p'
$ws.Range("V22").ClearFormats()

$ws.Range("W22").Value = 'ln alpha/p'
$ws.Range("W22").NumberFormat = '0.000'

$ws.Range("X22").Value = 'p/E'
$ws.Range("X22").NumberFormat = '0.000'

$ws.Range("T23").Value = 8000
$ws.Range("T23").NumberFormat = '0.000'

$ws.Range("U23").Formula = '=U11*100'
$ws.Range("U23").ClearFormats()

$ws.Range("V23").Formula = '=U23/$L$1'
$ws.Range("V23").ClearFormats()

$ws.Range("W23").Formula = '=LN(V23)'
$ws.Range("W23").ClearFormats()

$ws.Range("X23").Formula = '=$L$1/T23'
$ws.Range("X23").ClearFormats()

$ws.Range("T24").Value = 9000
$ws.Range("T24").NumberFormat = '0.000'

$ws.Range("U24").Formula = '=U13*100'
$ws.Range("U24").ClearFormats()

$ws.Range("V24").Formula = '=U24/$L$1'
$ws.Range("V24").ClearFormats()

$ws.Range("W24").Formula = '=LN(V24)'
$ws.Range("W24").ClearFormats()

$ws.Range("X24").Formula = '=$L$1/T24'
$ws.Range("X24").ClearFormats()

$ws.Range("T25").Value = 10000
$ws.Range("T25").NumberFormat = '0.000'

$ws.Range("U25").Formula = '=U15*100'
$ws.Range("U25").ClearFormats()

$ws.Range("V25").Formula = '=U25/$L$1'
$ws.Range("V25").ClearFormats()

$ws.Range("W25").Formula = '=LN(V25)'
$ws.Range("W25").ClearFormats()

$ws.Range("X25").Formula = '=$L$1/T25'
$ws.Range("X25").ClearFormats()

$ws.Range("T26").Value = 11000
$ws.Range("T26").NumberFormat = '0.000'

$ws.Range("U26").Formula = '=U17*100'
$ws.Range("U26").ClearFormats()

$ws.Range("V26").Formula = '=U26/$L$1'
$ws.Range("V26").ClearFormats()

$ws.Range("W26").Formula = '=LN(V26)'
$ws.Range("W26").ClearFormats()

$ws.Range("X26").Formula = '=$L$1/T26'
$ws.Range("X26").ClearFormats()

$ws.Range("T27").Value = 12000
$ws.Range("T27").NumberFormat = '0.000'

$ws.Range("U27").Formula = '=U19*100'
$ws.Range("U27").ClearFormats()

$ws.Range("V27").Formula = '=U27/$L$1'
$ws.Range("V27").ClearFormats()

$ws.Range("W27").Formula = '=LN(V27)'
$ws.Range("W27").ClearFormats()

$ws.Range("X27").Formula = '=$L$1/T27'
$ws.Range("X27").ClearFormats()

# Column widths for the new columns (T and W)
$ws.Columns.Item(20).ColumnWidth = 9.5703125
$ws.Columns.Item(23).ColumnWidth = 12.140625

# Selection moved as part of the edit session
$ws.Range("M30").Select()
